# "more updates to the ftau model"
# - Add Status Date (column E) entries for rows 17, 22 and 23
# - Append a new task row (24): "read all again about first order systems "
# - Update the saved view (scrolled down to the new row, selection on E16)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New Status Date values on existing rows - copy the date style (numFmtId 14)
# from the neighbouring Ideation Date cell so no new style/numFmt is minted.
$ws.Range("C17").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = (Get-Date -Year 2019 -Month 6 -Day 20).Date

$ws.Range("C22").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = (Get-Date -Year 2019 -Month 6 -Day 20).Date

$ws.Range("C23").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = (Get-Date -Year 2019 -Month 6 -Day 20).Date

# New row 24 - "Task Name" (wrap text style like B22/B23), "Ideation Date"
# (date style like C23) and "Status" = ToDo
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").Value = "read all again about first order systems "

$ws.Range("C23").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Range("C24").Value = (Get-Date -Year 2019 -Month 6 -Day 21).Date

$ws.Range("F24").Value = "ToDo"

# Update the saved selection/cursor state
$ws.Range("E16").Select() | Out-Null
